$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix "collecteur/" -> "collecteurs/" path typo (rows 2 and 3 only,
#     the other rows already use the "collecteurs/" prefix) ---
$ws.Range("A2").Value = "collecteurs/collecteur-2-rouge.png"
$ws.Range("A3").Value = "collecteurs/collecteur-3-rouge.png"

# --- Fix taille (size) values ---
# Row 18 "Console simple pr coll. 4/4 (paire)" : taille cleared (was bogus "4/4")
$ws.Range("C18").Value = ""

# Row 20 "Collecteur 4/4 x2 sorties rouges" : taille 4/4-3 -> 4/4-2 (matches the x2 row 19)
$ws.Range("C20").Value = "4/4-2"

# Row 21 "Collecteur 4/4 x5 sorties rouges" : taille 4/4-4 -> 4/4-5 (matches the x5 row 22)
$ws.Range("C21").Value = "4/4-5"

# --- Restore the view: no frozen/scrolled top-left cell, selection on C18 ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C18").Select()
